$d = $word.ActiveDocument

# The document currently ends with:
#   ... "Freundin will zurück zur Mama, Zettel am Kühlschrank" [+ _GoBack bookmark]
#   <empty paragraph>
# We need to insert a new "SW9" joke block between the "Freundin..." paragraph
# and the trailing empty paragraph, and move the _GoBack bookmark to the very
# end of the newly added content.

# Remove the existing _GoBack bookmark; it will be re-created at the new end.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# The final (trailing) empty paragraph of the document. It stays exactly as
# it is (empty, plain, de-CH) and now serves as the separator line right
# after the "Freundin..." joke, just like in the target document.
$trailing = $d.Paragraphs($d.Paragraphs.Count)

# Create three new, clean paragraph shells right before the trailing paragraph:
#   pHeader -> "SW9:" (bold)
#   pBody1  -> first body line of the joke
#   pBody2  -> punchline of the joke
$pHeader = $d.Paragraphs.Add($trailing.Range)
$pBody1  = $d.Paragraphs.Add($trailing.Range)
$pBody2  = $d.Paragraphs.Add($trailing.Range)

# Fill in the header paragraph.
$pHeader.Range.InsertAfter("SW9:")
$pHeader.Range.Font.Bold = $true
$pHeader.Range.LanguageID = "de-CH"

# Fill in the first body paragraph.
$pBody1.Range.InsertAfter("PM und junger Ing im Zug. Grossmutter mit hübscher Enkelin. Beleuchtung defekt, viele Tunnels, ")
$pBody1.Range.LanguageID = "de-CH"

# Fill in the punchline paragraph.
$pBody2.Range.InsertAfter("Der junge Ingenieur denkt: Das Leben ist schön, selten kann man eine junge, attraktive Frau küssen und gleichzeitig seinem PM eine Ohrfeige verpassen")
$pBody2.Range.LanguageID = "de-CH"

# Re-create the _GoBack bookmark at the very end of the new content.
$endRange = $pBody2.Range
$endRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $endRange)
